$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The paragraph "Nous affichons ... putProperty. " originally has
#    the word "putProperty" wrapped in a pair of proofErr spell-check
#    markers, splitting the sentence across three runs. Re-typing /
#    editing at the end of that paragraph (where the cursor/_GoBack
#    bookmark sits) makes Word drop the stale proofErr markers and
#    re-merge the adjoining identically-formatted runs. We reproduce
#    that via a Find/Replace of the whole sentence with itself, which
#    forces the engine to normalize the run layout.
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Nous affichons tout au long*putProperty*") {
        $targetPara = $cand
        break
    }
}

$needle = "Nous affichons tout au long du déroulé de l’algorithme le compteur de voisin de chaque nœud grâce à la méthode putProperty. "
$targetPara.Range.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null

# ------------------------------------------------------------------
# 2. The cursor was left right after "putProperty. " (that's where the
#    _GoBack bookmark lives, collapsed, at the very end of the
#    paragraph's text and before its paragraph mark). The user pressed
#    Enter, Enter, typed "TP_4 :", then Enter - which pushes the
#    bookmark into its own trailing paragraph.
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$pos = $bm.Range.Start

# Split the bookmark off into its own paragraph.
$d.Range($pos, $pos).InsertBefore("`r")

# Insert the "TP_4 :" paragraph right before that (now separate)
# bookmark paragraph.
$d.Range($pos + 1, $pos + 1).InsertBefore("TP_4 :`r")

# Insert a genuinely blank paragraph (no run) right before "TP_4 :",
# matching the self-closed <w:p/> produced by Word for an untouched
# empty paragraph.
$blankFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$d.Range($pos + 1, $pos + 1).InsertXML($blankFrag) | Out-Null

"done"
